$d = $word.ActiveDocument

# 1. Fecha: 16-02-2024 -> 31-10-2023
$d.Content.Find.Execute("16-02-2024", $true, $false, $false, $false, $false, $true, 1, $false, "31-10-2023", 2) | Out-Null

# 2. Nombre de la empresa: Julian Cándido Espinosa Trinidad -> prueba ti
$d.Content.Find.Execute("Julian Cándido Espinosa Trinidad", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 3. Estado: Ciudad de México -> cprueba ti
$d.Content.Find.Execute("Ciudad de México", $true, $false, $false, $false, $false, $true, 1, $false, "cprueba ti", 2) | Out-Null

# 4. Contacto: Angela Espinosa Trinidad -> prueba ti
$d.Content.Find.Execute("Angela Espinosa Trinidad", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 5. Correo: aaron.cuevas@splittel.com -> ramon.olea@splittel.com
$d.Content.Find.Execute("aaron.cuevas@splittel.com", $true, $false, $false, $false, $false, $true, 1, $false, "ramon.olea@splittel.com", 2) | Out-Null

# 6. Teléfono: 4421917076 -> 34543545
$d.Content.Find.Execute("4421917076", $true, $false, $false, $false, $false, $true, 1, $false, "34543545", 2) | Out-Null

# 7. Servicio: Calibración -> Reparación (standalone one-cell table, Find would be ambiguous
#    with the "Calibración (C)" inside the longer sentence, so target the table cell directly)
$tServicio = $d.Tables.Item(4)
$tServicio.Cell(1, 1).Range.Text = "Reparación"

# 8. Marca: Wandel & goltermann -> prueba ti
$d.Content.Find.Execute("Wandel & goltermann", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 9. Modelo: OLP_15B -> prueba ti
$d.Content.Find.Execute("OLP_15B", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 10. Número de serie: BG-0156 -> 4ggfgfg
$d.Content.Find.Execute("BG-0156", $true, $false, $false, $false, $false, $true, 1, $false, "4ggfgfg", 2) | Out-Null

# 11. Observaciones: aaaaa -> prueba ti
$d.Content.Find.Execute("aaaaa", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 12-14. Cantidad / No. Serie / Descripción first data row: empty -> 1 / dsfsd / prueba ti
$tDatos = $d.Tables.Item(6)
$tDatos.Cell(2, 1).Range.Text = "1"
$tDatos.Cell(2, 2).Range.Text = "dsfsd"
$tDatos.Cell(2, 3).Range.Text = "prueba ti"

# 15. Paquetería: DHL -> prueba ti
$d.Content.Find.Execute("DHL", $true, $false, $false, $false, $false, $true, 1, $false, "prueba ti", 2) | Out-Null

# 16. Numero de guía: 4680956733 -> dsfdsf59295
$d.Content.Find.Execute("4680956733", $true, $false, $false, $false, $false, $true, 1, $false, "dsfdsf59295", 2) | Out-Null
